# Applies the "Add files via upload" edit:
#  - Sheet "DANH SACH NO" (sheet index 1): fill in rows 20 and 21 (debtor #19
#    "Nguyen Huynh Anh Thu" / Nap quan huy / 20000, and debtor #20
#    "Huynh Quoc Phu" / Nap so / 75000) with their formulas/dates/status.
#  - Sheet "THONG KE NAP " (sheet index 2): log the same two payments as new
#    rows 135/136.
#  - Active sheet moves from "THONG KE NAP " back to "DANH SACH NO", with a
#    new selection on each sheet.

$wb = $excel.ActiveWorkbook
$wsDebt = $wb.Worksheets.Item(1)
$wsLog  = $wb.Worksheets.Item(2)

# ---- Sheet "DANH SACH NO": row 20 (A20 = 19) ----
$wsDebt.Range("B20").Value = "Nguyễn Huỳnh Anh Thư"
$wsDebt.Range("C20").Value = "Nạp quân huy"
$wsDebt.Range("D20").Value = 20000
$wsDebt.Range("E20").Value = 0
$wsDebt.Range("G20").Value = 0
$wsDebt.Range("H20").Value = 0
$wsDebt.Range("J20").Value = 46024
$wsDebt.Range("K20").Value = 46030
$wsDebt.Range("M20").Value = "Chưa trả đủ"

# ---- Sheet "DANH SACH NO": row 21 (A21 = 20) ----
$wsDebt.Range("B21").Value = "Huỳnh Quốc Phú"
$wsDebt.Range("C21").Value = "Nạp sò"
$wsDebt.Range("D21").Value = 75000
$wsDebt.Range("E21").Value = 0
$wsDebt.Range("G21").Value = 0
$wsDebt.Range("H21").Value = 0
$wsDebt.Range("J21").Value = 46024
$wsDebt.Range("K21").Value = 46030
$wsDebt.Range("M21").Value = "Chưa trả đủ"

# F20:F21 and I20:I21 are filled together as shared-formula groups, just like
# row 19's pattern extended down.
$wsDebt.Range("F20:F21").Formula = "=(D20+I20)-E20"
$wsDebt.Range("I20:I21").Formula = "=D20*H20"

# ---- Sheet "THONG KE NAP ": rows 135/136 log the two new payments ----
$wsLog.Range("A135").Value = 46024
$wsLog.Range("B135").Value = "Nguyễn Huỳnh Anh Thư"
$wsLog.Range("C135").Value = 20000
$wsLog.Range("D135").Value = "Nạp quân huy"

$wsLog.Range("A136").Value = 46024
$wsLog.Range("B136").Value = "Huỳnh Quốc Phú"
$wsLog.Range("C136").Value = 75000
$wsLog.Range("D136").Value = "Nạp sò"

# ---- View state: active sheet goes back to "DANH SACH NO" ----
[void]$wsLog.Range("D137").Select()
[void]$wsDebt.Activate()
[void]$wsDebt.Range("M26").Select()
